$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("branch")
$arr = ,@("90.0")
$ws.Range("S2").Value = $arr
